# Daily attendance processing - 2025-12-30 08:41:18
# Applies the recorded changes to the "Session Analysis Results" sheet:
#  1. Updates Class Statistics: Missing Sessions (L7) and Pending Sessions (L8)
#  2. Normalizes the "Recorded By" text order for sessions recorded by both
#     the instructor and the System (swap to "System, <email>")
#  3. Updates Group Statistics P/Q counters (rows 21-26) for the B1D1..B1F2
#     groups now that the 30/12/2025 session is no longer "Pending"
#  4. Flips the 30/12/2025 sessions (one per group) from "Pending" to
#     "Not Recorded" now that their window has elapsed unrecorded, including
#     re-coloring the row to match the existing "Not Recorded" style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Class Statistics summary box (K2:L10)
$ws.Range("L7").Value = 9
$ws.Range("L8").Value = 120

# 2. Swap "Recorded By" text order from "<email>, System" to "System, <email>"
$recordedByRows = @(8,9,10,12,14,15,17,34,35,36,38,40,41,43,60,61,62,64,66,67,69,86,87,88,90,92,93,95,112,113,114,116,118,119,121,138,139,140,142,144,145,147,164,167,170,191,194,197,218,221,224,245,248,251,272,275,278,299,302,305)
foreach ($r in $recordedByRows) {
    $cell = $ws.Range("G" + $r)
    $cell.Value = "System, dnasr281@gmail.com"
}

# 3. Group Statistics P/Q columns for rows 21-26 (B1D1, B1D2, B1E1, B1E2, B1F1, B1F2)
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 11

$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 11

$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 11

$ws.Range("P24").Value = 2
$ws.Range("Q24").Value = 11

$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 11

$ws.Range("P26").Value = 1
$ws.Range("Q26").Value = 11

# 4. Flip the 30/12/2025 "Pending" session rows to "Not Recorded" and restyle
#    the row to match the existing "Not Recorded" look (copy formats from a
#    known "Not Recorded" row, e.g. row 3).
$formatSource = $ws.Range("A3:I3")
$notRecordedRows = @(173,200,227,254,281,308)
foreach ($r in $notRecordedRows) {
    $destRow = $ws.Range("A" + $r + ":I" + $r)
    $formatSource.Copy()
    $destRow.PasteSpecial(-4122)
    $ws.Range("I" + $r).Value = "Not Recorded"
}

$excel.CutCopyMode = $false
